# Updated symbol list on Sat Dec 24 22:29:04 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Price" (column D) cells hold numeric-looking values that are actually
# stored as text. When writing them back through COM, Excel would otherwise
# auto-convert a numeric-looking string into a real number (losing the exact
# textual formatting, e.g. trailing zeros). To avoid that we briefly force a
# Text number format before assigning the value, then restore the cell style
# so the cell keeps looking like the rest of the (unformatted) column.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Simple price (column D) updates
Set-TextValue $ws.Range("D2")  "244.51"
Set-TextValue $ws.Range("D3")  "21.86"
Set-TextValue $ws.Range("D4")  "5.408"
Set-TextValue $ws.Range("D5")  "0.06040"
Set-TextValue $ws.Range("D6")  "3.395"
Set-TextValue $ws.Range("D7")  "0.8138"
Set-TextValue $ws.Range("D8")  "0.9268"
Set-TextValue $ws.Range("D9")  "0.1440"
Set-TextValue $ws.Range("D10") "0.07433"
Set-TextValue $ws.Range("D11") "0.03378"
Set-TextValue $ws.Range("D12") "0.03051"
Set-TextValue $ws.Range("D13") "0.09429"
Set-TextValue $ws.Range("D15") "0.001589"
Set-TextValue $ws.Range("D17") "0.0005943"
Set-TextValue $ws.Range("D18") "0.005683"
Set-TextValue $ws.Range("D19") "0.004155"
Set-TextValue $ws.Range("D20") "0.0009889"
Set-TextValue $ws.Range("D21") "3.653"
Set-TextValue $ws.Range("D22") "6.434"
Set-TextValue $ws.Range("D26") "0.00008506"
Set-TextValue $ws.Range("D27") "0.0002901"
Set-TextValue $ws.Range("D40") "0.04009"

# Rows 41-43: the coin list shifted by one position (a new coin was inserted
# ahead of them on the live site), so the Coin/Link/Price/Volume columns
# rotate between these three rows, each also picking up a refreshed price.
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.006414"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1074"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.002902"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue $ws.Range("D44") "0.006393"
Set-TextValue $ws.Range("D45") "0.00005237"
Set-TextValue $ws.Range("D48") "0.002320"
